$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.692.56"
$ws.Range("E2").Value = "  +3.14%  "
$ws.Range("D3").Value = "3.147.25"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.95"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.55"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.88%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.147.34"
$ws.Range("E8").Value = "  +2.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.51"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("D15").Value = "68.647.49"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").Value = "3.672.57"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.59%  "
$ws.Range("D19").Value = "3.147.90"
$ws.Range("E19").Value = "  +2.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.48"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.86"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.702"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("E25").Value = "  +7.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.07"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.62"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.90%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.12"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.35%  "
$ws.Range("E30").Value = "  +5.56%  "
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.78"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.55"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.962"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("E39").Value = "  +8.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.06"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.15%  "
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.23"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.44"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.75"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "402.25"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.06"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +14.38%  "
$ws.Range("D47").Value = "2.811.07"
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "134.93"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("E51").Value = "  +11.20%  "
